$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the empty "ListParagraph" paragraph that carries the old
#    "_GoBack" bookmark, right before the "Agent Parameters:" heading.
#    Deleting the whole paragraph (incl. its end mark) removes the bookmark
#    and lets Word renumber the remaining bookmark ids automatically.
# ---------------------------------------------------------------------------
$i = 0
$targetIndex = -1
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($targetIndex -eq -1 -and $p.Range.Text.Trim() -eq "") {
        $nextIndex = $i + 1
        if ($nextIndex -le $d.Paragraphs.Count) {
            $nextText = $d.Paragraphs.Item($nextIndex).Range.Text.Trim()
            if ($nextText -eq "Agent Parameters:") {
                $targetIndex = $i
            }
        }
    }
}
if ($targetIndex -ne -1) {
    $d.Paragraphs.Item($targetIndex).Range.Delete()
}

# ---------------------------------------------------------------------------
# 2. LR_ACTOR / LR_CRITIC: bump the learning rate from 1e-3 to 1e-4
#    (the extra space that used to pad the comment column is consumed by
#    the edit, exactly as in the tracked change).
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("LR_ACTOR = 1e-3 ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "LR_ACTOR = 1e-4", 2) | Out-Null

$r2 = $d.Content
$r2.Find.Execute("LR_CRITIC = 1e-3 ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "LR_CRITIC = 1e-4", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Relocate "_GoBack" into the "model gets" sentence, splitting it into
#    "model get" + "s a +30 score ...". The stale bookmark was already
#    removed in step 1, so this simply (re)creates it at the new spot.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$found = $r3.Find.Execute("model get", $false, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
if ($found) {
    $r3.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $r3)
}
